$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lower AHP")

$ws.Range("K1").Value = 0.5
$ws.Range("A2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("K2").Value = 0.5
$ws.Range("L2").Value = 0.5
$ws.Range("A3").Value = 0.5
$ws.Range("B3").Value = 1
$ws.Range("F3").Value = 0.5
$ws.Range("K3").Value = 0.5
$ws.Range("L3").Value = 1
$ws.Range("A4").Value = 0.2
$ws.Range("B4").Value = 0.5
$ws.Range("C4").Value = 0.5
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.25
$ws.Range("G4").Value = 0.5
$ws.Range("H4").Value = 0.3333333333333333
$ws.Range("I4").Value = 0.25
$ws.Range("J4").Value = 0.25
$ws.Range("K4").Value = 0.2
$ws.Range("L4").Value = 0.5
$ws.Range("A5").Value = 0.25
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 0.5
$ws.Range("I5").Value = 0.3333333333333333
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("K5").Value = 0.25
$ws.Range("L5").Value = 0.5
$ws.Range("A6").Value = 1
$ws.Range("K6").Value = 1
$ws.Range("A7").Value = 0.5
$ws.Range("B7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("H7").Value = 1
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 0.25
$ws.Range("L7").Value = 0.5
$ws.Range("A8").Value = 0.3333333333333333
$ws.Range("B8").Value = 0.5
$ws.Range("C8").Value = 0.5
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("I8").Value = 0.5
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 0.25
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("A9").Value = 0.3333333333333333
$ws.Range("B9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("K9").Value = 0.3333333333333333
$ws.Range("L9").Value = 0.5
$ws.Range("A10").Value = 0.3333333333333333
$ws.Range("B10").Value = 0.5
$ws.Range("C10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("I10").Value = 1
$ws.Range("K10").Value = 0.2
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("F12").Value = 1
